# Add a new worksheet named "12" as the last tab, mirroring the header-row
# pattern used by the other sheets (Section Number / value / Course Color /
# R / G / B), with the RGB triple 18, 164, 136. In every existing sheet all
# six header values -- including the numeric-looking ones -- are stored as
# text, so force Text format on the numeric-looking cells before writing
# them (otherwise Excel auto-coerces a digit-only string to a number).

$wb = $excel.ActiveWorkbook

$previouslyActiveSheet = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Insert the new sheet after the current last tab ("fasf") so it lands at
# the end, instead of Excel's default of inserting before the active sheet.
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "12"

$newSheet.Range("A1").Value = "Secttion Number:"

$newSheet.Range("B1").NumberFormat = "@"
$newSheet.Range("B1").Value = "12"

$newSheet.Range("C1").Value = "Course Color: "

$newSheet.Range("D1").NumberFormat = "@"
$newSheet.Range("D1").Value = "18"

$newSheet.Range("E1").NumberFormat = "@"
$newSheet.Range("E1").Value = "164"

$newSheet.Range("F1").NumberFormat = "@"
$newSheet.Range("F1").Value = "136"

# Adding a sheet makes it active; restore the original active tab.
$previouslyActiveSheet.Activate()
